# Updates crypto price/volume data per the Fri Mar 22 04:41:47 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.913.73"
$ws.Range("E2").Value = "  -1.06%  "

$ws.Range("D3").Value = "3.506.44"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").Value = "'578.32"
$ws.Range("E5").Value = "  +5.21%  "

$ws.Range("D6").Value = "'178.10"
$ws.Range("E6").Value = "  -5.74%  "

$ws.Range("D7").Value = "'0.634"
$ws.Range("E7").Value = "  +4.71%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.20%  "

$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("D10").Value = "'0.156"
$ws.Range("E10").Value = "  +4.35%  "

$ws.Range("D11").Value = "'55.10"
$ws.Range("E11").Value = "  -0.11%  "

$ws.Range("E12").Value = "  +1.66%  "

$ws.Range("D13").Value = "'9.21"
$ws.Range("E13").Value = "  -1.56%  "

$ws.Range("D14").Value = "4.071.07"
$ws.Range("E14").Value = "  +0.00%  "

$ws.Range("D15").Value = "3.508.90"
$ws.Range("E15").Value = "  -0.10%  "

$ws.Range("E16").Value = "  +0.20%  "

$ws.Range("E17").Value = "  +0.89%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "65.869.34"
$ws.Range("E18").Value = "  -1.36%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'12.02"
$ws.Range("E19").Value = "  +2.25%  "

$ws.Range("E20").Value = "  +1.31%  "

$ws.Range("D21").Value = "'413.77"
$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("D22").Value = "'4.25"
$ws.Range("E22").Value = "  +8.65%  "

$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").Value = "'4.34"
$ws.Range("E23").Value = "  +3.09%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'85.85"
$ws.Range("E24").Value = "  +0.79%  "

$ws.Range("D25").Value = "'13.29"
$ws.Range("E25").Value = "  +12.15%  "

$ws.Range("E26").Value = "  -1.56%  "

$ws.Range("D27").Value = "'2.85"
$ws.Range("E27").Value = "  -2.49%  "

$ws.Range("D28").Value = "'9.06"
$ws.Range("E28").Value = "  +2.86%  "

$ws.Range("D29").Value = "'30.41"
$ws.Range("E29").Value = "  +0.71%  "

$ws.Range("D30").Value = "'623.23"
$ws.Range("E30").Value = "  -4.67%  "

$ws.Range("E31").Value = "  -3.24%  "

$ws.Range("D32").Value = "'11.65"
$ws.Range("E32").Value = "  -0.24%  "

$ws.Range("E33").Value = "  -0.63%  "

$ws.Range("E34").Value = "  +14.97%  "

$ws.Range("D35").Value = "'59.67"
$ws.Range("E35").Value = "  +0.31%  "

$ws.Range("E36").Value = "  -0.19%  "

$ws.Range("D37").Value = "0.0₃0795"
$ws.Range("E37").Value = "  -1.80%  "

$ws.Range("D38").Value = "'37.21"
$ws.Range("E38").Value = "  -3.81%  "

$ws.Range("D39").Value = "'3.54"
$ws.Range("E39").Value = "  +5.56%  "

$ws.Range("D40").Value = "3.293.42"
$ws.Range("E40").Value = "  +10.06%  "

$ws.Range("E41").Value = "  -2.98%  "

$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("D43").Value = "'2.90"
$ws.Range("E43").Value = "  -0.17%  "

$ws.Range("D44").Value = "'0.0416"
$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'3.26"
$ws.Range("E45").Value = "  -2.46%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.50"
$ws.Range("E46").Value = "  -5.03%  "

$ws.Range("E47").Value = "  -0.04%  "

$ws.Range("D48").Value = "'0.132"
$ws.Range("E48").Value = "  +1.88%  "

$ws.Range("D49").Value = "'140.07"
$ws.Range("E49").Value = "  +1.05%  "

$ws.Range("D50").Value = "'8.55"
$ws.Range("E50").Value = "  -4.09%  "

$ws.Range("D51").Value = "'2.32"
$ws.Range("E51").Value = "  -3.87%  "
